# feat: add 2022-Q4 data
#
# Before:  Sheet1=总计 (summary), Sheet2=2022-Q3 (fund detail), Sheet3=2022-Q2 (fund detail)
# After:   Sheet1=总计 (summary, +1 row), Sheet2=2022-Q4 (NEW fund detail),
#          Sheet3=2022-Q3 (old fund detail, unchanged), Sheet4=2022-Q2 (old fund detail, unchanged)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Duplicate the current "2022-Q3" sheet so the old Q3 numbers survive
#    untouched on a sheet of their own; the original sheet object becomes the
#    new "2022-Q4" sheet (it keeps its position right after 总计).
# ---------------------------------------------------------------------------
$sQ3 = $wb.Worksheets.Item(2)
$sQ3.Copy($null, $sQ3)
$sQ3Copy = $wb.Worksheets.Item(3)

$sQ3.Name = "2022-Q4"
$sQ3Copy.Name = "2022-Q3"

# ---------------------------------------------------------------------------
# 2) Update the fund figures on the (renamed) "2022-Q4" sheet. These columns
#    are stored as text in the workbook, so force text formatting before
#    assigning, then restore the "Normal" style so no stray number-format /
#    quote-prefix style gets attached to the cell.
# ---------------------------------------------------------------------------
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $sQ3.Range("D2") "1.27"
Set-TextValue $sQ3.Range("E2") "94.90"
Set-TextValue $sQ3.Range("F2") "3.95"
Set-TextValue $sQ3.Range("G2") "0.0502"

Set-TextValue $sQ3.Range("D3") "0.63"
Set-TextValue $sQ3.Range("E3") "94.90"
Set-TextValue $sQ3.Range("F3") "3.95"
Set-TextValue $sQ3.Range("G3") "0.0249"

# ---------------------------------------------------------------------------
# 3) Update the "总计" summary sheet: shift the existing two rows down one
#    slot (quarter labels move down), add a brand-new row for "2022-Q4" at
#    the top, and append the row that used to be missing (old row 3 -> row 4).
# ---------------------------------------------------------------------------
$sTotal = $wb.Worksheets.Item(1)

# Create row 4 with the same formatting as row 3 (style s="2" on column A).
$sTotal.Range("A3").Copy()
$sTotal.Range("A4").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

$sTotal.Range("A4").Value = 2
$sTotal.Range("B4").Value = "2022-Q2"
$sTotal.Range("C4").Value = 2
$sTotal.Range("D4").Value = 0.11

$sTotal.Range("B3").Value = "2022-Q3"
$sTotal.Range("D3").Value = 0.08

$sTotal.Range("B2").Value = "2022-Q4"

# ---------------------------------------------------------------------------
# 4) Restore the originally-active tab. Copying a sheet makes the new copy
#    active; the workbook originally had "2022-Q2" selected, so reselect it
#    (it is now the 4th tab) to keep that state unchanged.
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("2022-Q2").Activate()
